# Update Icam1-Il2ra LR-pair sheet with recomputed TPM-derived values.
# Ligand avg/total expression (G,H) depend on the Sending cluster (col A);
# Receptor avg/total expression (M,N) depend on the Target cluster (col D);
# the specificity + edge-weight columns (I,J,O,P,Q,R,S,T) are recomputed from those.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 30.87085333333333
$ws.Range("H2").Value = 92.61256
$ws.Range("I2").Value = 0.2985789950947061
$ws.Range("J2").Value = 0.2985789950947061
$ws.Range("M2").Value = 0.2196916666666667
$ws.Range("N2").Value = 0.6590750000000001
$ws.Range("O2").Value = 0.1091447240078814
$ws.Range("P2").Value = 0.1091447240078814
$ws.Range("Q2").Value = 6.782069220222223
$ws.Range("R2").Value = 61.03862298200001
$ws.Range("S2").Value = 0.03258832201416227
$ws.Range("T2").Value = 0.03258832201416227

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 30.87085333333333
$ws.Range("H3").Value = 92.61256
$ws.Range("I3").Value = 0.2985789950947061
$ws.Range("J3").Value = 0.2985789950947061
$ws.Range("O3").Value = 0.02865791113152811
$ws.Range("P3").Value = 0.02865791113152811
$ws.Range("Q3").Value = 1.78075430368
$ws.Range("R3").Value = 16.02678873312
$ws.Range("S3").Value = 0.008556650307165055
$ws.Range("T3").Value = 0.008556650307165056

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 30.87085333333333
$ws.Range("H4").Value = 92.61256
$ws.Range("I4").Value = 0.2985789950947061
$ws.Range("J4").Value = 0.2985789950947061
$ws.Range("M4").Value = 1.200477666666667
$ws.Range("N4").Value = 3.601433
$ws.Range("O4").Value = 0.5964077090132023
$ws.Range("P4").Value = 0.5964077090132023
$ws.Range("Q4").Value = 37.05976997760889
$ws.Range("R4").Value = 333.53792979848
$ws.Range("S4").Value = 0.1780748144238978
$ws.Range("T4").Value = 0.1780748144238979

# Row 5: ECs -> Resolving-Mac
$ws.Range("G5").Value = 30.87085333333333
$ws.Range("H5").Value = 92.61256
$ws.Range("I5").Value = 0.2985789950947061
$ws.Range("J5").Value = 0.2985789950947061
$ws.Range("M5").Value = 0.534994
$ws.Range("N5").Value = 1.604982
$ws.Range("O5").Value = 0.2657896558473883
$ws.Range("P5").Value = 0.2657896558473883
$ws.Range("Q5").Value = 16.51572130821333
$ws.Range("R5").Value = 148.64149177392
$ws.Range("S5").Value = 0.07935920834948099
$ws.Range("T5").Value = 0.079359208349481

# Row 6: FAPs -> ECs
$ws.Range("G6").Value = 33.793597
$ws.Range("I6").Value = 0.3268474027571036
$ws.Range("J6").Value = 0.3268474027571037
$ws.Range("M6").Value = 0.2196916666666667
$ws.Range("N6").Value = 0.6590750000000001
$ws.Range("O6").Value = 0.1091447240078814
$ws.Range("P6").Value = 0.1091447240078814
$ws.Range("Q6").Value = 7.424171647591668
$ws.Range("R6").Value = 66.81754482832501
$ws.Range("S6").Value = 0.03567366956661693
$ws.Range("T6").Value = 0.03567366956661693

# Row 7: FAPs -> FAPs
$ws.Range("G7").Value = 33.793597
$ws.Range("I7").Value = 0.3268474027571036
$ws.Range("J7").Value = 0.3268474027571037
$ws.Range("O7").Value = 0.02865791113152811
$ws.Range("P7").Value = 0.02865791113152811
$ws.Range("Q7").Value = 1.949349849348
$ws.Range("S7").Value = 0.009366763821783851
$ws.Range("T7").Value = 0.009366763821783853

# Row 8: FAPs -> MuSCs
$ws.Range("G8").Value = 33.793597
$ws.Range("I8").Value = 0.3268474027571036
$ws.Range("J8").Value = 0.3268474027571037
$ws.Range("M8").Value = 1.200477666666667
$ws.Range("N8").Value = 3.601433
$ws.Range("O8").Value = 0.5964077090132023
$ws.Range("P8").Value = 0.5964077090132023
$ws.Range("Q8").Value = 40.56845847483367
$ws.Range("R8").Value = 365.116126273503
$ws.Range("S8").Value = 0.1949343106752796
$ws.Range("T8").Value = 0.1949343106752796

# Row 9: FAPs -> Resolving-Mac
$ws.Range("G9").Value = 33.793597
$ws.Range("I9").Value = 0.3268474027571036
$ws.Range("J9").Value = 0.3268474027571037
$ws.Range("M9").Value = 0.534994
$ws.Range("N9").Value = 1.604982
$ws.Range("O9").Value = 0.2657896558473883
$ws.Range("P9").Value = 0.2657896558473883
$ws.Range("Q9").Value = 18.079371633418
$ws.Range("R9").Value = 162.714344700762
$ws.Range("S9").Value = 0.0868726586934233
$ws.Range("T9").Value = 0.08687265869342331

# Row 10: MuSCs -> ECs
$ws.Range("G10").Value = 2.981185666666667
$ws.Range("H10").Value = 8.943557
$ws.Range("I10").Value = 0.02883365130639111
$ws.Range("J10").Value = 0.02883365130639111
$ws.Range("M10").Value = 0.2196916666666667
$ws.Range("N10").Value = 0.6590750000000001
$ws.Range("O10").Value = 0.1091447240078814
$ws.Range("P10").Value = 0.1091447240078814
$ws.Range("Q10").Value = 0.6549416477527779
$ws.Range("R10").Value = 5.894474829775
$ws.Range("S10").Value = 0.003147040913975546
$ws.Range("T10").Value = 0.003147040913975546

# Row 11: MuSCs -> FAPs
$ws.Range("G11").Value = 2.981185666666667
$ws.Range("H11").Value = 8.943557
$ws.Range("I11").Value = 0.02883365130639111
$ws.Range("J11").Value = 0.02883365130639111
$ws.Range("O11").Value = 0.02865791113152811
$ws.Range("P11").Value = 0.02865791113152811
$ws.Range("Q11").Value = 0.171966713996
$ws.Range("R11").Value = 1.547700425964
$ws.Range("S11").Value = 0.0008263122167360258
$ws.Range("T11").Value = 0.0008263122167360258

# Row 12: MuSCs -> MuSCs
$ws.Range("G12").Value = 2.981185666666667
$ws.Range("H12").Value = 8.943557
$ws.Range("I12").Value = 0.02883365130639111
$ws.Range("J12").Value = 0.02883365130639111
$ws.Range("M12").Value = 1.200477666666667
$ws.Range("N12").Value = 3.601433
$ws.Range("O12").Value = 0.5964077090132023
$ws.Range("P12").Value = 0.5964077090132023
$ws.Range("Q12").Value = 3.578846813020112
$ws.Range("R12").Value = 32.209621317181
$ws.Range("S12").Value = 0.01719661191813025
$ws.Range("T12").Value = 0.01719661191813025

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("G13").Value = 2.981185666666667
$ws.Range("H13").Value = 8.943557
$ws.Range("I13").Value = 0.02883365130639111
$ws.Range("J13").Value = 0.02883365130639111
$ws.Range("M13").Value = 0.534994
$ws.Range("N13").Value = 1.604982
$ws.Range("O13").Value = 0.2657896558473883
$ws.Range("P13").Value = 0.2657896558473883
$ws.Range("Q13").Value = 1.594916444552667
$ws.Range("R13").Value = 14.354248000974
$ws.Range("S13").Value = 0.007663686257549291
$ws.Range("T13").Value = 0.007663686257549291

# Row 14: Resolving-Mac -> ECs
$ws.Range("G14").Value = 35.74694633333333
$ws.Range("H14").Value = 107.240839
$ws.Range("I14").Value = 0.3457399508417991
$ws.Range("J14").Value = 0.3457399508417991
$ws.Range("M14").Value = 0.2196916666666667
$ws.Range("N14").Value = 0.6590750000000001
$ws.Range("O14").Value = 0.1091447240078814
$ws.Range("P14").Value = 0.1091447240078814
$ws.Range("Q14").Value = 7.85330621821389
$ws.Range("R14").Value = 70.67975596392502
$ws.Range("S14").Value = 0.03773569151312665
$ws.Range("T14").Value = 0.03773569151312664

# Row 15: Resolving-Mac -> FAPs
$ws.Range("G15").Value = 35.74694633333333
$ws.Range("H15").Value = 107.240839
$ws.Range("I15").Value = 0.3457399508417991
$ws.Range("J15").Value = 0.3457399508417991
$ws.Range("O15").Value = 0.02865791113152811
$ws.Range("P15").Value = 0.02865791113152811
$ws.Range("Q15").Value = 2.062026852292
$ws.Range("R15").Value = 18.558241670628
$ws.Range("S15").Value = 0.009908184785843177
$ws.Range("T15").Value = 0.009908184785843177

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("G16").Value = 35.74694633333333
$ws.Range("H16").Value = 107.240839
$ws.Range("I16").Value = 0.3457399508417991
$ws.Range("J16").Value = 0.3457399508417991
$ws.Range("M16").Value = 1.200477666666667
$ws.Range("N16").Value = 3.601433
$ws.Range("O16").Value = 0.5964077090132023
$ws.Range("P16").Value = 0.5964077090132023
$ws.Range("Q16").Value = 42.91341072469856
$ws.Range("R16").Value = 386.220696522287
$ws.Range("S16").Value = 0.2062019719958946
$ws.Range("T16").Value = 0.2062019719958946

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("G17").Value = 35.74694633333333
$ws.Range("H17").Value = 107.240839
$ws.Range("I17").Value = 0.3457399508417991
$ws.Range("J17").Value = 0.3457399508417991
$ws.Range("M17").Value = 0.534994
$ws.Range("N17").Value = 1.604982
$ws.Range("O17").Value = 0.2657896558473883
$ws.Range("P17").Value = 0.2657896558473883
$ws.Range("Q17").Value = 19.12440180665533
$ws.Range("R17").Value = 172.119616259898
$ws.Range("S17").Value = 0.09189410254693474
$ws.Range("T17").Value = 0.09189410254693474
